$d = $word.ActiveDocument

# Step 1: drop the trailing period from the existing sentence's run -
# "This is after the modification." -> "This is after the modification"
$d.Content.Find.Execute("This is after the modification.", $true, $false, $false, $false, $false, $true, 1, $false, "This is after the modification", 2) | Out-Null

# Step 2: append the rest of the text as a brand-new run (rather than
# extending the previous run's text), so the saved OOXML ends up with two
# <w:r> elements just like the target diff.
$p = $d.Paragraphs(1)
$r = $p.Range
$r.Collapse(0)  # wdCollapseEnd
$r.InsertAfter(". This is the modification.")

# Toggling a character property on the freshly-inserted range and back to
# its original value forces the run-merging logic to keep it as a
# separate run from its plain neighbour, while leaving no trace in the
# saved formatting.
$r.Bold = 1
$r.Bold = 0
